# Add a new DDL test case row (ddl_020) that checks information_schema.tables
# no longer lists a table after it has been dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 21

$ws.Cells.Item($row, 1).Value = "ddl_020"
$ws.Cells.Item($row, 2).Value = "y"
$ws.Cells.Item($row, 3).Value = '删除表后验证是否从information_schema.tables中移除'
$ws.Cells.Item($row, 4).Value = "DDL"
$ws.Cells.Item($row, 5).Value = "drop"
$ws.Cells.Item($row, 6).Value = "schema56"

# Table_value_ref (column G) stays empty for this case, but still carries the
# sheet's usual text number format.
$ws.Cells.Item($row, 7).NumberFormat = "@"

$ws.Cells.Item($row, 8).Value = 'drop table $schema56'
$ws.Cells.Item($row, 9).Value = 'select `TABLE_CATALOG`,`TABLE_SCHEMA`,`TABLE_NAME`,`TABLE_TYPE`,`ENGINE`,`VERSION`,`ROW_FORMAT`,`TABLE_ROWS`,`AVG_ROW_LENGTH`,`DATA_LENGTH`,`MAX_DATA_LENGTH`,`INDEX_LENGTH`,`DATA_FREE`,`AUTO_INCREMENT`,`UPDATE_TIME`,`CHECK_TIME`,`TABLE_COLLATION`,`CHECKSUM`,`CREATE_OPTIONS`,`TABLE_COMMENT` from information_schema.tables where `TABLE_NAME`=''$schema56'' or `TABLE_SCHEMA` in (''MYSQL'', ''INFORMATION_SCHEMA'')'
$ws.Cells.Item($row, 10).Value = "src/test/resources/io.dingodb.test/testdata/cases/ddl/expectedresult/ddl_020.csv"
$ws.Cells.Item($row, 11).Value = "csv_containsAll"

# Restore the sheet's active selection to where the author left off editing.
$ws.Activate()
$ws.Range("K24").Select()
